$wb = $excel.ActiveWorkbook

# Rename Sheet1 to the Chinese title used in the published workbook.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Name = "九九乘法表"

# Drop the two unused, empty worksheets.
$excel.DisplayAlerts = $false
$null = $wb.Worksheets.Item("Sheet2").Delete()
$null = $wb.Worksheets.Item("Sheet3").Delete()
$excel.DisplayAlerts = $true

# Keep the print-area defined name pointing at the renamed sheet.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "_xlnm.Print_Area" -or $n.Name -like "*Print_Area*") {
        $n.RefersTo = "=九九乘法表!`$A`$1:`$J`$10"
    }
}

# Move the saved selection from C16 to B10.
$null = $ws1.Range("B10").Select()
